# Auto-generated edit script: refresh market-price derived columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Leve profit" tables.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 22
$ws.Range("I8").Value = 22
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 66
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 73
$ws.Range("N8").ClearContents()
# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 231.55556
$ws.Range("I12").Value = 262.2857
$ws.Range("K12").Value = 262.2857
$ws.Range("M12").Value = -92.28570000000002
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2495
$ws.Range("J40").Value = 2495
$ws.Range("L40").Value = 2495
$ws.Range("N40").Value = -2845
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 3339.25
$ws.Range("I98").Value = 3287.7
$ws.Range("J98").Value = 3597
$ws.Range("K98").Value = 3287.7
$ws.Range("L98").Value = 3597
$ws.Range("M98").Value = -1789.7
$ws.Range("N98").Value = -6593
# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 974.9091
$ws.Range("I111").Value = 772.4
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 2317.2
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = 749.8000000000002
$ws.Range("N111").Value = -15134
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 3339.25
$ws.Range("I122").Value = 3287.7
$ws.Range("J122").Value = 3597
$ws.Range("K122").Value = 9863.099999999999
$ws.Range("L122").Value = 10791
$ws.Range("M122").Value = -7413.099999999999
$ws.Range("N122").Value = -15691
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2914.7856
$ws.Range("I132").Value = 3031.8
$ws.Range("K132").Value = 9095.400000000001
$ws.Range("M132").Value = -6565.400000000001

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 16 (Leve Item ID 3775)
$ws.Range("H16").Value = 11048.75
$ws.Range("I16").Value = 12773.25
$ws.Range("J16").Value = 7599.75
$ws.Range("K16").Value = 12773.25
$ws.Range("L16").Value = 7599.75
$ws.Range("M16").Value = -12486.25
$ws.Range("N16").Value = -8173.75
# Row 25 (Leve Item ID 2471)
$ws.Range("H25").Value = 479.8
$ws.Range("I25").Value = 549.75
$ws.Range("K25").Value = 549.75
$ws.Range("M25").Value = -147.75
# Row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 24000
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19727
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 54 (Leve Item ID 2376)
$ws.Range("H54").Value = 60000
$ws.Range("J54").Value = 70000
$ws.Range("L54").Value = 70000
$ws.Range("N54").Value = -70968

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2 (Leve Item ID 1820)
$ws.Range("H2").Value = 303.625
$ws.Range("I2").Value = 57.5
$ws.Range("J2").Value = 549.75
$ws.Range("K2").Value = 57.5
$ws.Range("L2").Value = 549.75
$ws.Range("M2").Value = 55.5
$ws.Range("N2").Value = -775.75
# Row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 16061.625
$ws.Range("I41").Value = 4250
$ws.Range("J41").Value = 19998.834
$ws.Range("K41").Value = 4250
$ws.Range("L41").Value = 19998.834
$ws.Range("M41").Value = -3822
$ws.Range("N41").Value = -20854.834
# Row 50 (Leve Item ID 1862)
$ws.Range("H50").Value = 29750
$ws.Range("J50").Value = 29750
$ws.Range("L50").Value = 29750
$ws.Range("N50").Value = -31000
# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 30500.445
$ws.Range("I59").Value = 17502
$ws.Range("J59").Value = 34214.285
$ws.Range("K59").Value = 17502
$ws.Range("L59").Value = 34214.285
$ws.Range("M59").Value = -16357
$ws.Range("N59").Value = -36504.285
# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 24697.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 24697.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 24697.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -25719.5
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 37493.75
$ws.Range("I68").Value = 19950
$ws.Range("K68").Value = 19950
$ws.Range("M68").Value = -19201
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 37493.75
$ws.Range("I71").Value = 19950
$ws.Range("K71").Value = 59850
$ws.Range("M71").Value = -56106
# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 38361.727
$ws.Range("I74").Value = 36494
$ws.Range("J74").Value = 38548.5
$ws.Range("K74").Value = 36494
$ws.Range("L74").Value = 38548.5
$ws.Range("M74").Value = -35620
$ws.Range("N74").Value = -40296.5
# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 38361.727
$ws.Range("I77").Value = 36494
$ws.Range("J77").Value = 38548.5
$ws.Range("K77").Value = 109482
$ws.Range("L77").Value = 115645.5
$ws.Range("M77").Value = -105114
$ws.Range("N77").Value = -124381.5

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1699.6666
$ws.Range("I68").Value = 1599.6
$ws.Range("J68").Value = 1771.1428
$ws.Range("K68").Value = 4798.799999999999
$ws.Range("L68").Value = 5313.428400000001
$ws.Range("M68").Value = -3987.799999999999
$ws.Range("N68").Value = -6935.428400000001
# Row 70 (Leve Item ID 12867)
$ws.Range("H70").Value = 2777
$ws.Range("I70").Value = 2777
$ws.Range("K70").Value = 8331
$ws.Range("M70").Value = -8016
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1699.6666
$ws.Range("I71").Value = 1599.6
$ws.Range("J71").Value = 1771.1428
$ws.Range("K71").Value = 14396.4
$ws.Range("L71").Value = 15940.2852
$ws.Range("M71").Value = -10340.4
$ws.Range("N71").Value = -24052.2852
# Row 73 (Leve Item ID 12867)
$ws.Range("H73").Value = 2777
$ws.Range("I73").Value = 2777
$ws.Range("K73").Value = 8331
$ws.Range("M73").Value = -7239
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1761.25
$ws.Range("J107").Value = 1760.25
$ws.Range("L107").Value = 5280.75
$ws.Range("N107").Value = -9120.75
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 1416.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1416.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4249.9998
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8589.9998
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 1336.2
$ws.Range("J122").Value = 1095.6
$ws.Range("L122").Value = 9860.4
$ws.Range("N122").Value = -14760.4
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2249.2666
$ws.Range("I131").Value = 1598.75
$ws.Range("J131").Value = 2485.818
$ws.Range("K131").Value = 4796.25
$ws.Range("L131").Value = 7457.454000000001
$ws.Range("M131").Value = 243.75
$ws.Range("N131").Value = -17537.454
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1993.4286
$ws.Range("J132").Value = 2665
$ws.Range("L132").Value = 23985
$ws.Range("N132").Value = -29045

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 21868.375
$ws.Range("I43").Value = 6390
$ws.Range("J43").Value = 31155.4
$ws.Range("K43").Value = 6390
$ws.Range("L43").Value = 31155.4
$ws.Range("M43").Value = -6239
$ws.Range("N43").Value = -31457.4
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3960
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = -5996
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3960
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = -29984
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1255.7646
$ws.Range("I102").Value = 882
$ws.Range("K102").Value = 882
$ws.Range("M102").Value = 740
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1953
$ws.Range("I113").Value = 1952.2222
$ws.Range("K113").Value = 1952.2222
$ws.Range("M113").Value = 217.7778000000001
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2500.5
$ws.Range("I132").Value = 2667.3333
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 8001.999899999999
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5471.999899999999
$ws.Range("N132").Value = -11060

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 56 (Leve Item ID 3668)
$ws.Range("H56").Value = 59999.5
$ws.Range("I56").Value = 59999.5
$ws.Range("K56").Value = 59999.5
$ws.Range("M56").Value = -59308.5
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 2305.8
$ws.Range("I61").Value = 2305.8
$ws.Range("K61").Value = 2305.8
$ws.Range("M61").Value = -2103.8
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 2305.8
$ws.Range("I113").Value = 2305.8
$ws.Range("K113").Value = 2305.8
$ws.Range("M113").Value = -135.8000000000002
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 9240.049999999999
$ws.Range("I132").Value = 9400.117
$ws.Range("K132").Value = 28200.351
$ws.Range("M132").Value = -25670.351

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16 (Leve Item ID 26304)
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 58 (Leve Item ID 3187)
$ws.Range("H58").Value = 24288.8
$ws.Range("I58").Value = 10450
$ws.Range("J58").Value = 45047
$ws.Range("K58").Value = 10450
$ws.Range("L58").Value = 45047
$ws.Range("M58").Value = -10142
$ws.Range("N58").Value = -45663
# Row 94 (Leve Item ID 18075)
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

